# PM&C supplied text updates to the "education_participation" workbook.
#
# Summary of changes:
#  1. Two footnotes in the Description sheet get trailing whitespace added
#     (as supplied by PM&C).
#  2. A new "Source" / citation row is appended to the Description sheet.
#  3. The new Source row (and the three footnote rows above it) use a
#     plain/normal wrapped-text style instead of the small italic note font,
#     and row heights grow to fit the (re-)wrapped text.
#  4. The previously-active "Data" tab is swapped for "Description" becoming
#     the active tab, with a new cell selection on each sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Description")

# --- 1. Footnote text updates (trailing spaces as supplied) ---------------
$ws2.Range("B7").Value = "People who were in school level education are excluded. "
$ws2.Range("B9").Value = "People permanently unable to work are excluded from the in-scope population.  "

# --- 2. New Source / citation row ------------------------------------------
$ws2.Range("A10").Value = "Source"
$ws2.Range("B10").Value = "ABS unpublished, 2006 Census of Population and Housing and 2011 Census of Population and Housing."

# --- 3. Re-style the footnote + source cells (plain font, wrapped) --------
# Derive the base "plain, wrapped" look once from an existing plain cell,
# then fan it out by copy/paste-format so every cell shares one style
# definition instead of each cell minting its own.
$ws2.Range("B3").Copy()
$ws2.Range("B7").PasteSpecial(-4122)
$ws2.Range("B7").WrapText = $true

$ws2.Range("B7").Copy()
$ws2.Range("B8").PasteSpecial(-4122)
$ws2.Range("B9").PasteSpecial(-4122)
$ws2.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Rows.Item(8).RowHeight = 49.45
$ws2.Rows.Item(9).RowHeight = 25.45
$ws2.Rows.Item(10).RowHeight = 25.45

# --- 4. Active sheet / selection changes -----------------------------------
[void]$ws1.Range("A22").Select()
[void]$ws2.Activate()
[void]$ws2.Range("B17").Select()
